$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking columns (Price, Hora) so exact formatting/precision is preserved
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '244.13'
$ws.Range('G2').Value = '18'

# Row 3
$ws.Range('D3').Value = '23.19'
$ws.Range('G3').Value = '18'

# Row 4
$ws.Range('D4').Value = '5.609'
$ws.Range('G4').Value = '18'

# Row 5
$ws.Range('D5').Value = '0.05864'
$ws.Range('G5').Value = '18'

# Row 6
$ws.Range('D6').Value = '3.408'
$ws.Range('G6').Value = '18'

# Row 7
$ws.Range('D7').Value = '6.479'
$ws.Range('G7').Value = '18'

# Row 8
$ws.Range('G8').Value = '18'

# Row 9
$ws.Range('D9').Value = '0.7982'
$ws.Range('G9').Value = '18'

# Row 10
$ws.Range('D10').Value = '0.1459'
$ws.Range('G10').Value = '18'

# Row 11
$ws.Range('D11').Value = '0.07629'
$ws.Range('G11').Value = '18'

# Row 12
$ws.Range('D12').Value = '0.03236'
$ws.Range('G12').Value = '18'

# Row 13
$ws.Range('D13').Value = '0.02971'
$ws.Range('G13').Value = '18'

# Row 14
$ws.Range('D14').Value = '0.09229'
$ws.Range('G14').Value = '18'

# Row 15
$ws.Range('D15').Value = '0.001664'
$ws.Range('G15').Value = '18'

# Row 16
$ws.Range('D16').Value = '3.344'
$ws.Range('G16').Value = '18'

# Row 17
$ws.Range('D17').Value = '0.04746'
$ws.Range('G17').Value = '18'

# Row 18
$ws.Range('D18').Value = '0.0005996'
$ws.Range('G18').Value = '18'

# Row 19
$ws.Range('D19').Value = '0.006242'
$ws.Range('G19').Value = '18'

# Row 20
$ws.Range('D20').Value = '0.001067'
$ws.Range('G20').Value = '18'

# Row 21
$ws.Range('D21').Value = '0.003830'
$ws.Range('G21').Value = '18'

# Row 22
$ws.Range('D22').Value = '0.0001499'
$ws.Range('G22').Value = '18'

# Row 23
$ws.Range('D23').Value = '3.691'
$ws.Range('G23').Value = '18'

# Row 24
$ws.Range('G24').Value = '18'

# Row 25
$ws.Range('D25').Value = '0.3333'
$ws.Range('G25').Value = '18'

# Row 26
$ws.Range('D26').Value = '0.1237'
$ws.Range('G26').Value = '18'

# Row 27
$ws.Range('D27').Value = '0.0004005'
$ws.Range('G27').Value = '18'

# Row 28
$ws.Range('G28').Value = '18'

# Row 29
$ws.Range('G29').Value = '18'

# Row 30
$ws.Range('G30').Value = '18'

# Row 31
$ws.Range('G31').Value = '18'

# Row 32
$ws.Range('G32').Value = '18'

# Row 33
$ws.Range('G33').Value = '18'

# Row 34
$ws.Range('G34').Value = '18'

# Row 35
$ws.Range('G35').Value = '18'

# Row 36
$ws.Range('G36').Value = '18'

# Row 37
$ws.Range('G37').Value = '18'

# Row 38
$ws.Range('G38').Value = '18'

# Row 39
$ws.Range('G39').Value = '18'

# Row 40
$ws.Range('D40').Value = '0.04308'
$ws.Range('G40').Value = '18'

# Row 41
$ws.Range('D41').Value = '0.007143'
$ws.Range('G41').Value = '18'

# Row 42
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = '0.1052'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('G42').Value = '18'

# Row 43
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = '0.003207'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('G43').Value = '18'

# Row 44
$ws.Range('D44').Value = '0.008805'
$ws.Range('G44').Value = '18'

# Row 45
$ws.Range('G45').Value = '18'

# Row 46
$ws.Range('D46').Value = '0.00005758'
$ws.Range('G46').Value = '18'

# Row 47
$ws.Range('G47').Value = '18'

# Row 48
$ws.Range('D48').Value = '0.7862'
$ws.Range('G48').Value = '18'

# Row 49
$ws.Range('D49').Value = '0.1052'
$ws.Range('G49').Value = '18'

# Row 50
$ws.Range('D50').Value = '0.00002103'
$ws.Range('G50').Value = '18'

# Row 51
$ws.Range('G51').Value = '18'
